$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.019.79"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "3.537.42"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.31"
$ws.Range("E5").Value = "  -2.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "196.73"
$ws.Range("E6").Value = "  +4.77%  "
$ws.Range("E7").Value = "  -1.77%  "
$ws.Range("E9").Value = "  -4.73%  "
$ws.Range("E10").Value = "  -1.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.02"
$ws.Range("E11").Value = "  +0.34%  "
$ws.Range("E12").Value = "  -2.03%  "
$ws.Range("E13").Value = "  -1.83%  "
$ws.Range("D14").Value = "4.100.61"
$ws.Range("E14").Value = "  -0.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "604.30"
$ws.Range("E15").Value = "  -1.98%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.23"
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").Value = "70.158.58"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("D19").Value = "3.532.99"
$ws.Range("E19").Value = "  -0.88%  "
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("E21").Value = "  -0.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.12"
$ws.Range("E22").Value = "  +2.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.31"
$ws.Range("E23").Value = "  +4.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "102.20"
$ws.Range("E24").Value = "  -2.79%  "
$ws.Range("E25").Value = "  -2.23%  "
$ws.Range("E26").Value = "  +2.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.93"
$ws.Range("E27").Value = "  -0.54%  "
$ws.Range("E28").Value = "  -2.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.57"
$ws.Range("E29").Value = "  -2.56%  "
$ws.Range("E30").Value = "  +0.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.33"
$ws.Range("E31").Value = "  +16.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.69"
$ws.Range("E32").Value = "  +1.18%  "
$ws.Range("E33").Value = "  -1.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.28"
$ws.Range("E34").Value = "  -1.54%  "
$ws.Range("D35").Value = "0.0₃0854"
$ws.Range("E35").Value = "  +8.80%  "
$ws.Range("D36").Value = "3.765.43"
$ws.Range("E36").Value = "  +6.14%  "
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.07"
$ws.Range("E38").Value = "  -3.17%  "
$ws.Range("E39").Value = "  +2.43%  "
$ws.Range("E40").Value = "  -1.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.69"
$ws.Range("E41").Value = "  -1.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "492.60"
$ws.Range("E42").Value = "  -8.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.133"
$ws.Range("E43").Value = "  -3.65%  "
$ws.Range("E44").Value = "  -2.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.84"
$ws.Range("E45").Value = "  -3.98%  "
$ws.Range("E46").Value = "  -2.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.32"
$ws.Range("E47").Value = "  -2.28%  "
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.66"
$ws.Range("E49").Value = "  -4.43%  "
$ws.Range("E50").Value = "  +2.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "130.50"
$ws.Range("E51").Value = "  -2.63%  "
